$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 67 ("Fuji royal" / 2020-11-26 batch).
# This pushes old rows 67-74 down to 71-78 unchanged, and leaves fresh blank
# rows 67-70 (inheriting column D's date-number-format style) ready to fill
# with the new week's data (2021-09-10, "Segunda" quality, $/caja 18 kilos granel).
$ws.Rows("67:70").Insert()

# Row 67: Fuji royal - Segunda
$ws.Cells.Item(67,1).Value = 1
$ws.Cells.Item(67,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(67,3).Value = "Arica y Parinacota"
$ws.Cells.Item(67,4).Value = 44449
$ws.Cells.Item(67,5).Value = 15
$ws.Cells.Item(67,6).Value = "Fruta"
$ws.Cells.Item(67,7).Value = 100104
$ws.Cells.Item(67,8).Value = "Frutos de pepita"
$ws.Cells.Item(67,9).Value = 100104002
$ws.Cells.Item(67,10).Value = "Manzana"
$ws.Cells.Item(67,11).Value = "Fuji royal"
$ws.Cells.Item(67,12).Value = "Segunda"
$ws.Cells.Item(67,13).Value = 250
$ws.Cells.Item(67,14).Value = 18000
$ws.Cells.Item(67,15).Value = 19000
$ws.Cells.Item(67,16).Value = 18500
$ws.Cells.Item(67,17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(67,18).Value = "Provincia de Curicó"
$ws.Cells.Item(67,19).Value = 1028
$ws.Cells.Item(67,20).Value = 18

# Row 68: Granny Smith - Segunda
$ws.Cells.Item(68,1).Value = 1
$ws.Cells.Item(68,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(68,3).Value = "Arica y Parinacota"
$ws.Cells.Item(68,4).Value = 44449
$ws.Cells.Item(68,5).Value = 15
$ws.Cells.Item(68,6).Value = "Fruta"
$ws.Cells.Item(68,7).Value = 100104
$ws.Cells.Item(68,8).Value = "Frutos de pepita"
$ws.Cells.Item(68,9).Value = 100104002
$ws.Cells.Item(68,10).Value = "Manzana"
$ws.Cells.Item(68,11).Value = "Granny Smith"
$ws.Cells.Item(68,12).Value = "Segunda"
$ws.Cells.Item(68,13).Value = 250
$ws.Cells.Item(68,14).Value = 18000
$ws.Cells.Item(68,15).Value = 19000
$ws.Cells.Item(68,16).Value = 18500
$ws.Cells.Item(68,17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(68,18).Value = "Provincia de Curicó"
$ws.Cells.Item(68,19).Value = 1028
$ws.Cells.Item(68,20).Value = 18

# Row 69: Royal Gala - Segunda
$ws.Cells.Item(69,1).Value = 1
$ws.Cells.Item(69,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(69,3).Value = "Arica y Parinacota"
$ws.Cells.Item(69,4).Value = 44449
$ws.Cells.Item(69,5).Value = 15
$ws.Cells.Item(69,6).Value = "Fruta"
$ws.Cells.Item(69,7).Value = 100104
$ws.Cells.Item(69,8).Value = "Frutos de pepita"
$ws.Cells.Item(69,9).Value = 100104002
$ws.Cells.Item(69,10).Value = "Manzana"
$ws.Cells.Item(69,11).Value = "Royal Gala"
$ws.Cells.Item(69,12).Value = "Segunda"
$ws.Cells.Item(69,13).Value = 270
$ws.Cells.Item(69,14).Value = 18000
$ws.Cells.Item(69,15).Value = 19000
$ws.Cells.Item(69,16).Value = 18500
$ws.Cells.Item(69,17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(69,18).Value = "Provincia de Curicó"
$ws.Cells.Item(69,19).Value = 1028
$ws.Cells.Item(69,20).Value = 18

# Row 70: Scarlett - Segunda
$ws.Cells.Item(70,1).Value = 1
$ws.Cells.Item(70,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(70,3).Value = "Arica y Parinacota"
$ws.Cells.Item(70,4).Value = 44449
$ws.Cells.Item(70,5).Value = 15
$ws.Cells.Item(70,6).Value = "Fruta"
$ws.Cells.Item(70,7).Value = 100104
$ws.Cells.Item(70,8).Value = "Frutos de pepita"
$ws.Cells.Item(70,9).Value = 100104002
$ws.Cells.Item(70,10).Value = "Manzana"
$ws.Cells.Item(70,11).Value = "Scarlett"
$ws.Cells.Item(70,12).Value = "Segunda"
$ws.Cells.Item(70,13).Value = 270
$ws.Cells.Item(70,14).Value = 18000
$ws.Cells.Item(70,15).Value = 19000
$ws.Cells.Item(70,16).Value = 18500
$ws.Cells.Item(70,17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(70,18).Value = "Provincia de Curicó"
$ws.Cells.Item(70,19).Value = 1028
$ws.Cells.Item(70,20).Value = 18
